# Automatically re-generate list and index
#
# The "Review date" column (B) stores dates as plain text (e.g. "2011-03-21").
# Every review date that falls on the 21st of its month is being rolled
# forward one day, to the 22nd. Dates that are not on the 21st (e.g. the
# 2024-06-02, 2024-06-26, 2024-10-24, ... one-off entries) are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRows = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $usedRows; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $current = [string]$cell.Text

    if ($current -match '^\d{4}-\d{2}-21$') {
        $updated = $current.Substring(0, $current.Length - 2) + "22"

        # Writing a "YYYY-MM-DD"-shaped string straight into .Value makes
        # Excel auto-coerce it into a date serial number. Force the cell to
        # Text format first so the literal string is preserved, then restore
        # the cell's original (default/"Normal") style so no visible
        # formatting change is introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $updated
        $cell.Style = "Normal"
    }
}
